# Automatic update of files.
# Refresh the "Förändrad" (changed) date stamp for every existing case to
# 2026-02-17 (serial 46070) and add the two newly published cases
# ("A 8679-2026" and "A 8929-2026"), which pushes the two oldest rows
# down to make room.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Bump the "Förändrad" column (C) to the new run date for every
#        row that already exists (rows 2-10) -----------------------------
for ($r = 2; $r -le 10; $r++) {
    $ws.Cells.Item($r, 3).Value = 46070
}

# --- 2. Re-key rows 6-10 with their new content. Row 6 ("A 50762-2025")
#        moves up one slot, two brand-new cases are inserted at rows 7-8,
#        and the previously-existing rows 6, 8, 9, 10 shift down to
#        9, 10, 11, 12 respectively. ------------------------------------

$rows = @(
    @{ Row = 6;  Case = "A 50762-2025"; Date = 45946;              Area = 2.7 },
    @{ Row = 7;  Case = "A 8679-2026";  Date = 46066;              Area = 2.1 },
    @{ Row = 8;  Case = "A 8929-2026";  Date = 46069.34543981482;  Area = 1.2 },
    @{ Row = 9;  Case = "A 14516-2023"; Date = 45012.86600694444;  Area = 0.4 },
    @{ Row = 10; Case = "A 23798-2024"; Date = 45455.43208333333;  Area = 1.3 },
    @{ Row = 11; Case = "A 4156-2023";  Date = 44953;              Area = 1.5 },
    @{ Row = 12; Case = "A 4159-2023";  Date = 44953;              Area = 0.5 }
)

foreach ($entry in $rows) {
    $r = $entry.Row

    $ws.Cells.Item($r, 1).Value = $entry.Case                 # A - Beteckning
    $ws.Cells.Item($r, 2).Value = $entry.Date                 # B - Datum
    $ws.Cells.Item($r, 2).NumberFormat = "YYYY-MM-DD"
    $ws.Cells.Item($r, 3).Value = 46070                       # C - Förändrad
    $ws.Cells.Item($r, 3).NumberFormat = "YYYY-MM-DD"
    $ws.Cells.Item($r, 4).Value = "SKÅNE LÄN"                 # D - Län
    $ws.Cells.Item($r, 5).Value = "ÅSTORP"                    # E - Kommun
    $ws.Cells.Item($r, 7).Value = $entry.Area                 # G - Area (ha)

    # H-Q: Fridlysta .. Alla arter, always 0 for these cases
    for ($col = 8; $col -le 17; $col++) {
        $ws.Cells.Item($r, $col).Value = 0
    }

    # R - Artnamn: always blank, wrap-text formatted
    $ws.Cells.Item($r, 18).WrapText = $true
}

# --- 3. Row heights: rows 2-9 already carry an explicit 15pt custom
#        height; row 10 gains one in this revision, as does the brand
#        new row 11. Row 12 (new) stays on the sheet default. -----------
$ws.Rows(10).RowHeight = 15
$ws.Rows(11).RowHeight = 15
